# Fixed albedo data script and created plot for it
$wb = $excel.ActiveWorkbook
$locations = $wb.Worksheets.Item("locations")
$albedos = $wb.Worksheets.Item("albedos")

# --- albedos sheet: fix the header label (drop stray leading space before "Albedo") ---
$albedos.Range("B1").Value = "Albedo"

# --- albedos sheet: the old Arial/bordered header+body styling is gone; reuse the
#     plain "Calibri 9.6, vertical-centered" style that the locations sheet already
#     uses for its timezone column ---
$locations.Range("D2").Copy()
$albedos.Range("A1:B10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# rows no longer need the taller custom height now that the bold header look is gone
$albedos.Rows.Item(1).RowHeight = 15
$albedos.Rows.Item(2).RowHeight = 15
$albedos.Rows.Item(8).RowHeight = 15
$albedos.Rows.Item(10).RowHeight = 15

# --- switch focus to the albedos sheet/tab, matching where the user ended up ---
$albedos.Select()
$albedos.Range("D6").Select()

$locations.Range("D4").Select()
$albedos.Activate()
$albedos.Range("D6").Select()
